$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.400.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.586.90'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.91%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.51%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.585.65'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.92%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("E10").Value = '  -0.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.93'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.17%  '

$ws.Range("E12").Value = '  +0.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.193.83'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000205'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.58'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.37%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.573.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.13%  '

$ws.Range("E17").Value = '  +1.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.429.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.05'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.83%  '

$ws.Range("E20").Value = '  +2.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '423.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.611'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.85%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000120'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.38'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.57%  '

$ws.Range("E28").Value = '  +3.80%  '

$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.582.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.86%  '

$ws.Range("E32").Value = '  +4.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.32%  '

$ws.Range("E34").Value = '  -1.84%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("E36").Value = '  -0.81%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.54%  '

$ws.Range("E38").Value = '  -2.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '175.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0853'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.880'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.68%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '46.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.04%  '

$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.22%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.32'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.87%  '

$ws.Range("E49").Value = '  +0.51%  '

$ws.Range("E50").Value = '  -4.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.942'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.26%  '
